$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = 15
$ws.Range("C1").Value = 16
$ws.Range("D1").Value = 15
$ws.Range("E1").Value = 16

$ws.Range("B2").Value = 1.9341576716069004
$ws.Range("D2").Value = 2.8616854811904409
$ws.Range("C2").ClearContents()
$ws.Range("E2").ClearContents()

$ws.Range("B3").Value = 1.6688906381123014
$ws.Range("C3").Value = -1.3172335172622707
$ws.Range("D3").Value = 3.1385564478826642
$ws.Range("E3").Value = -0.68186776113077985

$ws.Range("B1:E3").Select() | Out-Null
